$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.489.21"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "1.850.08"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.48%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5162"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3259"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06768"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7729"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.79%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.902.50"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07707"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.038"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007909"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.33%  "
$ws.Range("D20").Value = "26.515.68"
$ws.Range("D21").Value = "2.096.05"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.536"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.547"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.947"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.358"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.657"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.209"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.177"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08755"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04829"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.136"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.842"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6896"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01807"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.225"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4918"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "113.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9008"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.142"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.795"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4248"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.95%  "
$ws.Range("E47").Value = "  -7.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.112"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05901"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.38%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.09%  "
